# Update the "Förändrad" (changed) date column (C) for rows 2-16
# from 2023-11-03 (serial 45233) to 2023-11-13 (serial 45243).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = Get-Date -Year 2023 -Month 11 -Day 13 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
